$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 7.787422666666667
$ws.Cells.Item(2, 8).Value = 23.362268
$ws.Cells.Item(2, 9).Value = 0.2161047632645357
$ws.Cells.Item(2, 10).Value = 0.2161047632645357
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1528053333333333
$ws.Cells.Item(2, 14).Value = 0.458416
$ws.Cells.Item(2, 15).Value = 0.01103433215988526
$ws.Cells.Item(2, 16).Value = 0.01103433215988526
$ws.Cells.Item(2, 17).Value = 1.189959716387555
$ws.Cells.Item(2, 18).Value = 10.709637447488
$ws.Cells.Item(2, 19).Value = 0.002384571739194257
$ws.Cells.Item(2, 20).Value = 0.002384571739194257

$ws.Cells.Item(3, 7).Value = 7.787422666666667
$ws.Cells.Item(3, 8).Value = 23.362268
$ws.Cells.Item(3, 9).Value = 0.2161047632645357
$ws.Cells.Item(3, 10).Value = 0.2161047632645357
$ws.Cells.Item(3, 15).Value = 0.8539197603380489
$ws.Cells.Item(3, 16).Value = 0.8539197603380488
$ws.Cells.Item(3, 17).Value = 92.08804856570134
$ws.Cells.Item(3, 18).Value = 828.7924370913121
$ws.Cells.Item(3, 19).Value = 0.1845361276547631
$ws.Cells.Item(3, 20).Value = 0.1845361276547631

$ws.Cells.Item(4, 7).Value = 7.787422666666667
$ws.Cells.Item(4, 8).Value = 23.362268
$ws.Cells.Item(4, 9).Value = 0.2161047632645357
$ws.Cells.Item(4, 10).Value = 0.2161047632645357
$ws.Cells.Item(4, 13).Value = 1.712948333333333
$ws.Cells.Item(4, 14).Value = 5.138845
$ws.Cells.Item(4, 15).Value = 0.1236949029880405
$ws.Cells.Item(4, 16).Value = 0.1236949029880405
$ws.Cells.Item(4, 17).Value = 13.33945267782889
$ws.Cells.Item(4, 18).Value = 120.05507410046
$ws.Cells.Item(4, 19).Value = 0.0267310577272602
$ws.Cells.Item(4, 20).Value = 0.0267310577272602

$ws.Cells.Item(5, 7).Value = 7.787422666666667
$ws.Cells.Item(5, 8).Value = 23.362268
$ws.Cells.Item(5, 9).Value = 0.2161047632645357
$ws.Cells.Item(5, 10).Value = 0.2161047632645357
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1571906666666667
$ws.Cells.Item(5, 14).Value = 0.471572
$ws.Cells.Item(5, 15).Value = 0.01135100451402528
$ws.Cells.Item(5, 16).Value = 0.01135100451402528
$ws.Cells.Item(5, 17).Value = 1.224110160588445
$ws.Cells.Item(5, 18).Value = 11.016991445296
$ws.Cells.Item(5, 19).Value = 0.002453006143318109
$ws.Cells.Item(5, 20).Value = 0.002453006143318109

$ws.Cells.Item(6, 7).Value = 18.16892433333334
$ws.Cells.Item(6, 8).Value = 54.50677300000001
$ws.Cells.Item(6, 9).Value = 0.5041964793605992
$ws.Cells.Item(6, 10).Value = 0.5041964793605993
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1528053333333333
$ws.Cells.Item(6, 14).Value = 0.458416
$ws.Cells.Item(6, 15).Value = 0.01103433215988526
$ws.Cells.Item(6, 16).Value = 0.01103433215988526
$ws.Cells.Item(6, 17).Value = 2.776308539063111
$ws.Cells.Item(6, 18).Value = 24.986776851568
$ws.Cells.Item(6, 19).Value = 0.005563471427109584
$ws.Cells.Item(6, 20).Value = 0.005563471427109586

$ws.Cells.Item(7, 7).Value = 18.16892433333334
$ws.Cells.Item(7, 8).Value = 54.50677300000001
$ws.Cells.Item(7, 9).Value = 0.5041964793605992
$ws.Cells.Item(7, 10).Value = 0.5041964793605993
$ws.Cells.Item(7, 15).Value = 0.8539197603380489
$ws.Cells.Item(7, 16).Value = 0.8539197603380488
$ws.Cells.Item(7, 19).Value = 0.4305433368188909
$ws.Cells.Item(7, 20).Value = 0.430543336818891

$ws.Cells.Item(8, 7).Value = 18.16892433333334
$ws.Cells.Item(8, 8).Value = 54.50677300000001
$ws.Cells.Item(8, 9).Value = 0.5041964793605992
$ws.Cells.Item(8, 10).Value = 0.5041964793605993
$ws.Cells.Item(8, 13).Value = 1.712948333333333
$ws.Cells.Item(8, 14).Value = 5.138845
$ws.Cells.Item(8, 15).Value = 0.1236949029880405
$ws.Cells.Item(8, 16).Value = 0.1236949029880405
$ws.Cells.Item(8, 17).Value = 31.12242865524278
$ws.Cells.Item(8, 18).Value = 280.101857897185
$ws.Cells.Item(8, 19).Value = 0.06236653460142089
$ws.Cells.Item(8, 20).Value = 0.06236653460142089

$ws.Cells.Item(9, 7).Value = 18.16892433333334
$ws.Cells.Item(9, 8).Value = 54.50677300000001
$ws.Cells.Item(9, 9).Value = 0.5041964793605992
$ws.Cells.Item(9, 10).Value = 0.5041964793605993
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1571906666666667
$ws.Cells.Item(9, 14).Value = 0.471572
$ws.Cells.Item(9, 15).Value = 0.01135100451402528
$ws.Cells.Item(9, 16).Value = 0.01135100451402528
$ws.Cells.Item(9, 17).Value = 2.855985328572889
$ws.Cells.Item(9, 18).Value = 25.703867957156
$ws.Cells.Item(9, 19).Value = 0.005723136513177816
$ws.Cells.Item(9, 20).Value = 0.005723136513177816

$ws.Cells.Item(10, 7).Value = 6.195365666666667
$ws.Cells.Item(10, 8).Value = 18.586097
$ws.Cells.Item(10, 9).Value = 0.1719244078612872
$ws.Cells.Item(10, 10).Value = 0.1719244078612872
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.1528053333333333
$ws.Cells.Item(10, 14).Value = 0.458416
$ws.Cells.Item(10, 15).Value = 0.01103433215988526
$ws.Cells.Item(10, 16).Value = 0.01103433215988526
$ws.Cells.Item(10, 17).Value = 0.946684915816889
$ws.Cells.Item(10, 18).Value = 8.520164242352001
$ws.Cells.Item(10, 19).Value = 0.001897071022733031
$ws.Cells.Item(10, 20).Value = 0.001897071022733031

$ws.Cells.Item(11, 7).Value = 6.195365666666667
$ws.Cells.Item(11, 8).Value = 18.586097
$ws.Cells.Item(11, 9).Value = 0.1719244078612872
$ws.Cells.Item(11, 10).Value = 0.1719244078612872
$ws.Cells.Item(11, 15).Value = 0.8539197603380489
$ws.Cells.Item(11, 16).Value = 0.8539197603380488
$ws.Cells.Item(11, 17).Value = 73.26161155170534
$ws.Cells.Item(11, 18).Value = 659.3545039653482
$ws.Cells.Item(11, 19).Value = 0.1468096491571713
$ws.Cells.Item(11, 20).Value = 0.1468096491571713

$ws.Cells.Item(12, 7).Value = 6.195365666666667
$ws.Cells.Item(12, 8).Value = 18.586097
$ws.Cells.Item(12, 9).Value = 0.1719244078612872
$ws.Cells.Item(12, 10).Value = 0.1719244078612872
$ws.Cells.Item(12, 13).Value = 1.712948333333333
$ws.Cells.Item(12, 14).Value = 5.138845
$ws.Cells.Item(12, 15).Value = 0.1236949029880405
$ws.Cells.Item(12, 16).Value = 0.1236949029880405
$ws.Cells.Item(12, 17).Value = 10.61234129310722
$ws.Cells.Item(12, 18).Value = 95.511071637965
$ws.Cells.Item(12, 19).Value = 0.02126617295167822
$ws.Cells.Item(12, 20).Value = 0.02126617295167823

$ws.Cells.Item(13, 7).Value = 6.195365666666667
$ws.Cells.Item(13, 8).Value = 18.586097
$ws.Cells.Item(13, 9).Value = 0.1719244078612872
$ws.Cells.Item(13, 10).Value = 0.1719244078612872
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.1571906666666667
$ws.Cells.Item(13, 14).Value = 0.471572
$ws.Cells.Item(13, 15).Value = 0.01135100451402528
$ws.Cells.Item(13, 16).Value = 0.01135100451402528
$ws.Cells.Item(13, 17).Value = 0.9738536593871113
$ws.Cells.Item(13, 18).Value = 8.764682934484
$ws.Cells.Item(13, 19).Value = 0.001951514729704594
$ws.Cells.Item(13, 20).Value = 0.001951514729704594

$ws.Cells.Item(14, 7).Value = 3.883692333333334
$ws.Cells.Item(14, 8).Value = 11.651077
$ws.Cells.Item(14, 9).Value = 0.1077743495135779
$ws.Cells.Item(14, 10).Value = 0.1077743495135779
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.1528053333333333
$ws.Cells.Item(14, 14).Value = 0.458416
$ws.Cells.Item(14, 15).Value = 0.01103433215988526
$ws.Cells.Item(14, 16).Value = 0.01103433215988526
$ws.Cells.Item(14, 17).Value = 0.5934489015591111
$ws.Cells.Item(14, 18).Value = 5.341040114032
$ws.Cells.Item(14, 19).Value = 0.001189217970848387
$ws.Cells.Item(14, 20).Value = 0.001189217970848387

$ws.Cells.Item(15, 7).Value = 3.883692333333334
$ws.Cells.Item(15, 8).Value = 11.651077
$ws.Cells.Item(15, 9).Value = 0.1077743495135779
$ws.Cells.Item(15, 10).Value = 0.1077743495135779
$ws.Cells.Item(15, 15).Value = 0.8539197603380489
$ws.Cells.Item(15, 16).Value = 0.8539197603380488
$ws.Cells.Item(15, 17).Value = 45.92554732351868
$ws.Cells.Item(15, 18).Value = 413.3299259116681
$ws.Cells.Item(15, 19).Value = 0.09203064670722358
$ws.Cells.Item(15, 20).Value = 0.09203064670722357

$ws.Cells.Item(16, 7).Value = 3.883692333333334
$ws.Cells.Item(16, 8).Value = 11.651077
$ws.Cells.Item(16, 9).Value = 0.1077743495135779
$ws.Cells.Item(16, 10).Value = 0.1077743495135779
$ws.Cells.Item(16, 13).Value = 1.712948333333333
$ws.Cells.Item(16, 14).Value = 5.138845
$ws.Cells.Item(16, 15).Value = 0.1236949029880405
$ws.Cells.Item(16, 16).Value = 0.1236949029880405
$ws.Cells.Item(16, 17).Value = 6.652564309562778
$ws.Cells.Item(16, 18).Value = 59.873078786065
$ws.Cells.Item(16, 19).Value = 0.01333113770768119
$ws.Cells.Item(16, 20).Value = 0.01333113770768119

$ws.Cells.Item(17, 7).Value = 3.883692333333334
$ws.Cells.Item(17, 8).Value = 11.651077
$ws.Cells.Item(17, 9).Value = 0.1077743495135779
$ws.Cells.Item(17, 10).Value = 0.1077743495135779
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.1571906666666667
$ws.Cells.Item(17, 14).Value = 0.471572
$ws.Cells.Item(17, 15).Value = 0.01135100451402528
$ws.Cells.Item(17, 16).Value = 0.01135100451402528
$ws.Cells.Item(17, 17).Value = 0.6104801870048889
$ws.Cells.Item(17, 18).Value = 5.494321683044
$ws.Cells.Item(17, 19).Value = 0.001223347127824761
$ws.Cells.Item(17, 20).Value = 0.001223347127824761

Write-Host "Applied $($wb.Name) updates"
